# Update "viewed/sold" counts (column F) across sheets to reflect the
# regenerated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3021
$ws1.Range("F6").Value = 22
$ws1.Range("F7").Value = 248
$ws1.Range("F10").Value = 7112
$ws1.Range("F11").Value = 51
$ws1.Range("F12").Value = 128
$ws1.Range("F13").Value = 400
$ws1.Range("F14").Value = 638
$ws1.Range("F16").Value = 2297
$ws1.Range("F17").Value = 1559
$ws1.Range("F18").Value = 150
$ws1.Range("F19").Value = 1150
$ws1.Range("F20").Value = 18
$ws1.Range("F21").Value = 214
$ws1.Range("F22").Value = 370
$ws1.Range("F23").Value = 69
$ws1.Range("F24").Value = 69
$ws1.Range("F25").Value = 1819
$ws1.Range("F26").Value = 1734
$ws1.Range("F27").Value = 1042
$ws1.Range("F29").Value = 1698
$ws1.Range("F30").Value = 1277
$ws1.Range("F34").Value = 1076
$ws1.Range("F35").Value = 460
$ws1.Range("F37").Value = 2533
$ws1.Range("F38").Value = 2817
$ws1.Range("F39").Value = 2095
$ws1.Range("F41").Value = 196
$ws1.Range("F44").Value = 26
$ws1.Range("F45").Value = 34
$ws1.Range("F47").Value = 131
$ws1.Range("F48").Value = 190
$ws1.Range("F49").Value = 46
$ws1.Range("F50").Value = 421

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F18").Value = 75

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 1746
$ws3.Range("F7").Value = 1870
$ws3.Range("F8").Value = 2809
$ws3.Range("F9").Value = 1069
$ws3.Range("F10").Value = 988
$ws3.Range("F12").Value = 355
$ws3.Range("F13").Value = 1681
$ws3.Range("F14").Value = 7732

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3021
$ws4.Range("F6").Value = 248
$ws4.Range("F7").Value = 1746
$ws4.Range("F9").Value = 2809
$ws4.Range("F10").Value = 7112
$ws4.Range("F11").Value = 1069
$ws4.Range("F12").Value = 988
$ws4.Range("F13").Value = 51
$ws4.Range("F14").Value = 400
$ws4.Range("F15").Value = 355
$ws4.Range("F16").Value = 638
$ws4.Range("F18").Value = 2297
$ws4.Range("F19").Value = 1559
$ws4.Range("F20").Value = 150
$ws4.Range("F21").Value = 1150
$ws4.Range("F22").Value = 18
$ws4.Range("F23").Value = 370
$ws4.Range("F25").Value = 69
$ws4.Range("F26").Value = 1819
$ws4.Range("F28").Value = 1042
$ws4.Range("F29").Value = 1698
$ws4.Range("F30").Value = 1277
$ws4.Range("F34").Value = 1076
$ws4.Range("F35").Value = 75
$ws4.Range("F37").Value = 460
$ws4.Range("F39").Value = 2533
$ws4.Range("F40").Value = 2817
$ws4.Range("F41").Value = 2095
$ws4.Range("F43").Value = 196
$ws4.Range("F47").Value = 190
$ws4.Range("F49").Value = 421
